$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.344.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").Value = "'2.592.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'316.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").Value = "'97.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +2.57%  "

$ws.Range("D10").Value = "'35.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").Value = "'7.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "'2.986.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.45%  "

$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").Value = "'2.595.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").Value = "'0.849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "'43.379.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("D19").Value = "'6.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.37%  "

$ws.Range("D20").Value = "'12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("D21").Value = "'0.0₃0967"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "'69.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'255.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.99%  "

$ws.Range("D24").Value = "'2.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("E25").Value = "  +3.55%  "

$ws.Range("D26").Value = "'27.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  +1.24%  "

$ws.Range("D29").Value = "'40.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("D30").Value = "'10.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("D31").Value = "'5.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("D32").Value = "'157.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("E33").Value = "  +4.81%  "

$ws.Range("D34").Value = "'2.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.55%  "

$ws.Range("D35").Value = "'0.0812"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.12%  "

$ws.Range("E36").Value = "  +3.54%  "

$ws.Range("D37").Value = "'18.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("E39").Value = "  +7.95%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "'22.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.52%  "

$ws.Range("D42").Value = "'4.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.07%  "

$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").Value = "'2.019.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").Value = "'9.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "

$ws.Range("D48").Value = "'83.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "

$ws.Range("D49").Value = "'2.837.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.42%  "

$ws.Range("D50").Value = "'75.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.37%  "

$ws.Range("E51").Value = "  +3.02%  "
